$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.769.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.73%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.67%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.27'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9986'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3631'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.79'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3263'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07068'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9972'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.041'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.58'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.659.06'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.625'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001049'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06623'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9975'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '79.00'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -5.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.915'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.76'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -8.31%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.708.70'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.393'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -13.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.99'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.63'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -8.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.843.20'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.220'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.60'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.079'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.845'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -13.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08440'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.665'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.30'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -9.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.278'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.222'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.06%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02239'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.92%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06039'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -9.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2070'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.205'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -9.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9983'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5926'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.819'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.73'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5642'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -7.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.25'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.951'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -7.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06974'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.192'
